# Macroferia Regional de Talca - Zanahoria: add a new weekly price record.
# A new data row is inserted at row 517 (pushing the existing rows 517-582
# down to 518-583) so the sheet grows from A1:R582 to A1:R583.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 517, shifting everything below it down by one.
$ws.Rows.Item(517).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A517").Value2 = 5
$ws.Range("B517").Value2 = "Macroferia Regional de Talca"
$ws.Range("C517").Value2 = "Maule"
$ws.Range("D517").Value2 = 45154
$ws.Range("E517").Value2 = 7
$ws.Range("F517").Value2 = 100114013
$ws.Range("G517").Value2 = "Zanahoria"
$ws.Range("H517").Value2 = "Sin especificar"
$ws.Range("I517").Value2 = "Primera"
$ws.Range("J517").Value2 = 700
$ws.Range("K517").Value2 = 5000
$ws.Range("L517").Value2 = 5000
$ws.Range("M517").Value2 = 5000
$ws.Range("N517").Value2 = "`$/saco 20 kilos"
$ws.Range("O517").Value2 = "Región de Ñuble"
$ws.Range("P517").Value2 = 250
$ws.Range("Q517").Value2 = 20
$ws.Range("R517").Value2 = "Hortaliza"
